$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '295.10'
Set-TextValue $ws.Range('E2') '-2.43%'
Set-TextValue $ws.Range('G2') '16'
Set-TextValue $ws.Range('D3') '31.97'
Set-TextValue $ws.Range('E3') '0.12%'
Set-TextValue $ws.Range('G3') '16'
Set-TextValue $ws.Range('D4') '5.009'
Set-TextValue $ws.Range('E4') '-2.08%'
Set-TextValue $ws.Range('G4') '16'
Set-TextValue $ws.Range('D5') '0.07450'
Set-TextValue $ws.Range('E5') '-4.37%'
Set-TextValue $ws.Range('G5') '16'
Set-TextValue $ws.Range('D6') '1.918'
Set-TextValue $ws.Range('E6') '-14.72%'
Set-TextValue $ws.Range('G6') '16'
Set-TextValue $ws.Range('D7') '7.775'
Set-TextValue $ws.Range('E7') '-0.71%'
Set-TextValue $ws.Range('G7') '16'
Set-TextValue $ws.Range('D8') '3.783'
Set-TextValue $ws.Range('E8') '-0.73%'
Set-TextValue $ws.Range('G8') '16'
Set-TextValue $ws.Range('D9') '0.9200'
Set-TextValue $ws.Range('E9') '-0.51%'
Set-TextValue $ws.Range('G9') '16'
Set-TextValue $ws.Range('D10') '0.1729'
Set-TextValue $ws.Range('E10') '-1.52%'
Set-TextValue $ws.Range('G10') '16'
Set-TextValue $ws.Range('D11') '0.07711'
Set-TextValue $ws.Range('E11') '-0.33%'
Set-TextValue $ws.Range('G11') '16'
Set-TextValue $ws.Range('D12') '0.08330'
Set-TextValue $ws.Range('E12') '-5.95%'
Set-TextValue $ws.Range('G12') '16'
Set-TextValue $ws.Range('D13') '0.03028'
Set-TextValue $ws.Range('E13') '-3.78%'
Set-TextValue $ws.Range('G13') '16'
Set-TextValue $ws.Range('D14') '0.09961'
Set-TextValue $ws.Range('E14') '-0.42%'
Set-TextValue $ws.Range('G14') '16'
Set-TextValue $ws.Range('D15') '0.001501'
Set-TextValue $ws.Range('E15') '-1.05%'
Set-TextValue $ws.Range('G15') '16'
Set-TextValue $ws.Range('D16') '0.005944'
Set-TextValue $ws.Range('E16') '-0.24%'
Set-TextValue $ws.Range('G16') '16'
Set-TextValue $ws.Range('G17') '16'
Set-TextValue $ws.Range('D18') '3.463'
Set-TextValue $ws.Range('E18') '0.65%'
Set-TextValue $ws.Range('G18') '16'
Set-TextValue $ws.Range('D19') '2.140'
Set-TextValue $ws.Range('E19') '-5.35%'
Set-TextValue $ws.Range('G19') '16'
Set-TextValue $ws.Range('D20') '0.3340'
Set-TextValue $ws.Range('E20') '2.07%'
Set-TextValue $ws.Range('G20') '16'
Set-TextValue $ws.Range('D21') '0.1333'
Set-TextValue $ws.Range('E21') '0.27%'
Set-TextValue $ws.Range('G21') '16'
Set-TextValue $ws.Range('D22') '4.401'
Set-TextValue $ws.Range('E22') '2.82%'
Set-TextValue $ws.Range('G22') '16'
Set-TextValue $ws.Range('D23') '0.1990'
Set-TextValue $ws.Range('E23') '9.37%'
Set-TextValue $ws.Range('G23') '16'
Set-TextValue $ws.Range('D24') '0.04541'
Set-TextValue $ws.Range('E24') '-0.90%'
Set-TextValue $ws.Range('G24') '16'
Set-TextValue $ws.Range('D25') '0.001240'
Set-TextValue $ws.Range('E25') '-0.49%'
Set-TextValue $ws.Range('G25') '16'
Set-TextValue $ws.Range('D26') '0.004061'
Set-TextValue $ws.Range('E26') '-9.17%'
Set-TextValue $ws.Range('G26') '16'
Set-TextValue $ws.Range('D27') '0.0001260'
Set-TextValue $ws.Range('E27') '0.55%'
Set-TextValue $ws.Range('G27') '16'
Set-TextValue $ws.Range('G28') '16'
Set-TextValue $ws.Range('G29') '16'
Set-TextValue $ws.Range('G30') '16'
Set-TextValue $ws.Range('G31') '16'
Set-TextValue $ws.Range('G32') '16'
Set-TextValue $ws.Range('G33') '16'
Set-TextValue $ws.Range('G34') '16'
Set-TextValue $ws.Range('G35') '16'
Set-TextValue $ws.Range('G36') '16'
Set-TextValue $ws.Range('G37') '16'
Set-TextValue $ws.Range('G38') '16'
Set-TextValue $ws.Range('D39') '0.01646'
Set-TextValue $ws.Range('E39') '-6.76%'
Set-TextValue $ws.Range('G39') '16'
Set-TextValue $ws.Range('D40') '0.04520'
Set-TextValue $ws.Range('E40') '-5.92%'
Set-TextValue $ws.Range('G40') '16'
Set-TextValue $ws.Range('D41') '0.007389'
Set-TextValue $ws.Range('E41') '2.28%'
Set-TextValue $ws.Range('G41') '16'
Set-TextValue $ws.Range('D42') '0.1327'
Set-TextValue $ws.Range('E42') '-2.90%'
Set-TextValue $ws.Range('G42') '16'
Set-TextValue $ws.Range('D43') '0.002251'
Set-TextValue $ws.Range('E43') '5.92%'
Set-TextValue $ws.Range('G43') '16'
Set-TextValue $ws.Range('D44') '0.01007'
Set-TextValue $ws.Range('E44') '0.75%'
Set-TextValue $ws.Range('G44') '16'
Set-TextValue $ws.Range('D45') '0.00006118'
Set-TextValue $ws.Range('E45') '-1.87%'
Set-TextValue $ws.Range('G45') '16'
Set-TextValue $ws.Range('D46') '0.00000000756'
Set-TextValue $ws.Range('E46') '0.57%'
Set-TextValue $ws.Range('G46') '16'
Set-TextValue $ws.Range('D47') '1.643'
Set-TextValue $ws.Range('E47') '104.25%'
Set-TextValue $ws.Range('G47') '16'
Set-TextValue $ws.Range('D48') '0.003024'
Set-TextValue $ws.Range('E48') '-15.35%'
Set-TextValue $ws.Range('G48') '16'
Set-TextValue $ws.Range('D49') '0.00002117'
Set-TextValue $ws.Range('E49') '0.57%'
Set-TextValue $ws.Range('G49') '16'
Set-TextValue $ws.Range('D50') '0.0002016'
Set-TextValue $ws.Range('E50') '0.57%'
Set-TextValue $ws.Range('G50') '16'
Set-TextValue $ws.Range('G51') '16'
